# ============================================================
# LOQ4076.xlsx edit script
# Rewrites the "Objetivos / Programa resumido / Programa /
# Metodo / Criterio / Norma de recuperacao / Bibliografia" block:
# three new long-text entries are inserted into the shared-strings
# table (learning objectives, short + full syllabus, assessment
# rubric, bibliography), which pushes every following row down by
# one and adds a brand-new last row (24). Columns A & B share one
# <col min="1" max="2"> entry in the original file; the edit also
# splits that into separate per-column <col> entries.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Clear cells that exist today but must be empty afterwards ----
$ws.Range("A13").Clear()
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# ---- 2. Write the (possibly new) text into every live cell ----
#         Plain assignment keeps each existing cell's style (s="1/2/3");
#         brand-new cells get their style copied in step 3 below.
$t = @'
Ementa atual:
'@
$ws.Range("B1").Value = $t
$t = @'
Ementa modificada (dados modificados em vermelho):
'@
$ws.Range("C1").Value = $t

$t = @'
LOQ4076
'@
$ws.Range("B2").Value = $t
$t = @'
LOQ4076
'@
$ws.Range("C2").Value = $t

$t = @'
Nome:
'@
$ws.Range("A3").Value = $t
$t = @'
 Termodinâmica Aplicada
'@
$ws.Range("B3").Value = $t
$t = @'
 Termodinâmica Aplicada
'@
$ws.Range("C3").Value = $t

$t = @'
Name:
'@
$ws.Range("A4").Value = $t
$t = @'
Applied Thermodynamics
'@
$ws.Range("B4").Value = $t
$t = @'
Applied Thermodynamics
'@
$ws.Range("C4").Value = $t

$t = @'
Créditos-aula:
'@
$ws.Range("A5").Value = $t
$t = @'
4
'@
$ws.Range("B5").Value = $t
$t = @'
4
'@
$ws.Range("C5").Value = $t

$t = @'
Créditos-trabalho
'@
$ws.Range("A6").Value = $t
$t = @'
0
'@
$ws.Range("B6").Value = $t
$t = @'
0
'@
$ws.Range("C6").Value = $t

$t = @'
Carga horária:
'@
$ws.Range("A7").Value = $t
$t = @'
60 h
'@
$ws.Range("B7").Value = $t
$t = @'
60 h
'@
$ws.Range("C7").Value = $t

$t = @'
Ativação:
'@
$ws.Range("A8").Value = $t
$t = @'
01/01/2012
'@
$ws.Range("B8").Value = $t
$t = @'
01/01/2012
'@
$ws.Range("C8").Value = $t

$t = @'
Semestre ideal:
'@
$ws.Range("A9").Value = $t
$t = @'
EA-4,EP-6
'@
$ws.Range("B9").Value = $t
$t = @'
EA-4,EP-6
'@
$ws.Range("C9").Value = $t

$t = @'
Objetivos:
'@
$ws.Range("A10").Value = $t
$t = @'
Ao final do curso os estudantes deverão: - Compreender os aspectos mássicos, energéticos e entrópicos, envolvendo sistemas termodinâmicos abertos e fechados; - Dominar e ser capaz de fazer predições básicas de propriedades termodinâmicas, usando equações cúbicas de estado e relações termodinâmicas; - Desenvolver uma metodologia para poder solucionar os problemas de engenharia, nos aspectos termodinâmicos; Dominar o uso de tabelas de propriedades termodinâmicas;
'@
$ws.Range("B10").Value = $t
$t = @'
Ao final do curso os estudantes deverão: - Compreender os aspectos mássicos, energéticos e entrópicos, envolvendo sistemas termodinâmicos abertos e fechados; - Dominar e ser capaz de fazer predições básicas de propriedades termodinâmicas, usando equações cúbicas de estado e relações termodinâmicas; - Desenvolver uma metodologia para poder solucionar os problemas de engenharia, nos aspectos termodinâmicos; Dominar o uso de tabelas de propriedades termodinâmicas;
'@
$ws.Range("C10").Value = $t

$t = @'
Objectives:
'@
$ws.Range("A11").Value = $t

$t = @'
Docentes responsáveis:
'@
$ws.Range("A12").Value = $t

$t = @'
8554681 - Pedro Felipe Arce Castillo
'@
$ws.Range("B13").Value = $t
$t = @'
8554681 - Pedro Felipe Arce Castillo
'@
$ws.Range("C13").Value = $t

$t = @'
Programa resumido:
'@
$ws.Range("A14").Value = $t
$t = @'
A primeira Lei da Termodinâmica. Efeitos Térmicos. A segunda lei da Termodinâmica. Propriedades termodinâmicas dos fluidos. Termodinâmica de processos de escoamento. Produção de potencia a partir de calor. Refrigeração e liquefação
'@
$ws.Range("B14").Value = $t
$t = @'
A primeira Lei da Termodinâmica. Efeitos Térmicos. A segunda lei da Termodinâmica. Propriedades termodinâmicas dos fluidos. Termodinâmica de processos de escoamento. Produção de potencia a partir de calor. Refrigeração e liquefação
'@
$ws.Range("C14").Value = $t

$t = @'
Short syllabus:
'@
$ws.Range("A15").Value = $t

$t = @'
Programa:
'@
$ws.Range("A16").Value = $t
$t = @'
1  A primeira Lei da Termodinâmica
     1.1- Energia interna
     1.2- Estado termodinâmico e funções de estado
     1.3- Entalpia
     1.4- A regra das fases
     1.5- O processo reversível
     1.6- Processos a volume constante e a pressão constante
     1.7- Capacidade calorífica
2  Efeitos Térmicos 
     2.1- Calores Latentes de Substâncias Puras.
     2.2- Calor de Reação Padrão
     2.3- Calor Padrão de Formação
     2.4- Calor Padrão de Combustão
     2.5- O processo reversível
     2.6- A variação da entalpia com a Temperatura
3- A segunda lei da Termodinâmica
    3.1- Enunciados da lei
    3.2- Máquinas térmicas
    3.3- Escalas de temperaturas termodinâmicas
    3.4- Entropia 
    3.5- Variações da entropia de um gás ideal
    3.6- A terceira lei da termodinâmica
4- Produção de potencia a partir de calor
   4.1- A planta de potencia a vapor (maquina a vapor)
   4.2- Motores de combustão interna
   4.3- O motor Otto
   4.4- O motor Diesel
   4.5- A planta de potencia com turbina a gás
5- Refrigeração e liquefação
    5.1- O refrigerador de Carnot
    5.2- O ciclo com compresso a vapor
    5.3- Comparação de ciclos de refrigeração
    5.4- Refrigeração por absorção
    5.5- A bomba a calor
    6.6- Processos de liquefação
6- Termodinâmica de soluções
    6.1- Relações fundamentais entre propriedades
    6.2- O potencial químico 
    6.3- Fugacidade e coeficiente de fugacidade
    6.4- A solução Ideal
    6.5- Modelos para a energia de Gibbs
    6.6- Propriedades de mistura
    6.7- Efeitos térmicos em processos de mistura
7- Equilíbrio de fases
    7.1- Equilíbrio e estabilidade
    7.2- Equilíbrio líquido-líquido
    7.3- Equilíbrio líquido-líquido-vapor
    7.4- Equilíbrio sólido-líquido
    7.5- Equilíbrio sólido-vapor
    7.6- Equilíbrio na adsorção de gases em sólidos
'@
$ws.Range("B16").Value = $t
$t = @'
1  A primeira Lei da Termodinâmica
     1.1- Energia interna
     1.2- Estado termodinâmico e funções de estado
     1.3- Entalpia
     1.4- A regra das fases
     1.5- O processo reversível
     1.6- Processos a volume constante e a pressão constante
     1.7- Capacidade calorífica
2  Efeitos Térmicos 
     2.1- Calores Latentes de Substâncias Puras.
     2.2- Calor de Reação Padrão
     2.3- Calor Padrão de Formação
     2.4- Calor Padrão de Combustão
     2.5- O processo reversível
     2.6- A variação da entalpia com a Temperatura
3- A segunda lei da Termodinâmica
    3.1- Enunciados da lei
    3.2- Máquinas térmicas
    3.3- Escalas de temperaturas termodinâmicas
    3.4- Entropia 
    3.5- Variações da entropia de um gás ideal
    3.6- A terceira lei da termodinâmica
4- Produção de potencia a partir de calor
   4.1- A planta de potencia a vapor (maquina a vapor)
   4.2- Motores de combustão interna
   4.3- O motor Otto
   4.4- O motor Diesel
   4.5- A planta de potencia com turbina a gás
5- Refrigeração e liquefação
    5.1- O refrigerador de Carnot
    5.2- O ciclo com compresso a vapor
    5.3- Comparação de ciclos de refrigeração
    5.4- Refrigeração por absorção
    5.5- A bomba a calor
    6.6- Processos de liquefação
6- Termodinâmica de soluções
    6.1- Relações fundamentais entre propriedades
    6.2- O potencial químico 
    6.3- Fugacidade e coeficiente de fugacidade
    6.4- A solução Ideal
    6.5- Modelos para a energia de Gibbs
    6.6- Propriedades de mistura
    6.7- Efeitos térmicos em processos de mistura
7- Equilíbrio de fases
    7.1- Equilíbrio e estabilidade
    7.2- Equilíbrio líquido-líquido
    7.3- Equilíbrio líquido-líquido-vapor
    7.4- Equilíbrio sólido-líquido
    7.5- Equilíbrio sólido-vapor
    7.6- Equilíbrio na adsorção de gases em sólidos
'@
$ws.Range("C16").Value = $t

$t = @'
Syllabus:
'@
$ws.Range("A17").Value = $t

$t = @'
Avaliação:
'@
$ws.Range("A18").Value = $t

$t = @'
Método:
'@
$ws.Range("A19").Value = $t
$t = @'
2 provas escritas
'@
$ws.Range("B19").Value = $t
$t = @'
2 provas escritas
'@
$ws.Range("C19").Value = $t

$t = @'
Critério:
'@
$ws.Range("A20").Value = $t
$t = @'
serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas.
'@
$ws.Range("B20").Value = $t
$t = @'
serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. A média da disciplina será a média aritmética das duas provas.
'@
$ws.Range("C20").Value = $t

$t = @'
Norma de recuperação:
'@
$ws.Range("A21").Value = $t
$t = @'
prova escrita com conteúdo de todo o semestre
'@
$ws.Range("B21").Value = $t
$t = @'
prova escrita com conteúdo de todo o semestre
'@
$ws.Range("C21").Value = $t

$t = @'
Bibliografia:
'@
$ws.Range("A22").Value = $t
$t = @'
1)Smith, J.M.; Van Ness, H.C.; Abott, M. M.  Introdução à Termodinâmica da Engenharia Química. 7ª ed.  ISBN 978-85-216-1553-8, LTC editora, 2007.
2)Koretsky, M. D. Termodinâmica para Engenharia Química, 1ª ed.  ISBN 978-85-216-1530-9, LTC editora, 2007.
3)Terron, L. R. Termodinâmica Química Aplicada. 1ª ed.  ISBN 978-85-204-2082-9, Editora Manole Ltda, 2009.
4)Moran, M. J.; Shapiro, H. N. Princípios de Termodinâmica para Engenharia - 1ª ed.  ISBN 978-85-216-1689-4, LTC editora, 2009.
5)Van Wilen, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica  6ª Edição  2004
6)Sandler, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999
'@
$ws.Range("B22").Value = $t
$t = @'
1)Smith, J.M.; Van Ness, H.C.; Abott, M. M.  Introdução à Termodinâmica da Engenharia Química. 7ª ed.  ISBN 978-85-216-1553-8, LTC editora, 2007.
2)Koretsky, M. D. Termodinâmica para Engenharia Química, 1ª ed.  ISBN 978-85-216-1530-9, LTC editora, 2007.
3)Terron, L. R. Termodinâmica Química Aplicada. 1ª ed.  ISBN 978-85-204-2082-9, Editora Manole Ltda, 2009.
4)Moran, M. J.; Shapiro, H. N. Princípios de Termodinâmica para Engenharia - 1ª ed.  ISBN 978-85-216-1689-4, LTC editora, 2009.
5)Van Wilen, J. Sonntag, Richard. E. Fundamentos da Termodinâmica Clássica  6ª Edição  2004
6)Sandler, S. I., Chemical and Engineering Thermodynamics, 3rd ed., John Wiley & Sons, 1999
'@
$ws.Range("C22").Value = $t

$t = @'
Requisitos:
'@
$ws.Range("A23").Value = $t

$t = @'
LOQ4053 -  Balanços de Massa e Energia  (Requisito fraco)
'@
$ws.Range("B24").Value = $t
$t = @'
LOQ4053 -  Balanços de Massa e Energia  (Requisito fraco)
'@
$ws.Range("C24").Value = $t

# ---- 3. Brand-new cells (rows 14B/C, 16B/C, 22B/C, 23A, 24B/C did not
#         exist in the original sheet) need their number format /
#         font / alignment copied over from a same-column sibling ----
$ws.Range("B19").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("A19").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- 4. Row heights ----
# 4a. Rows whose explicit height changed (or are brand new)
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(24).RowHeight = 30

# 4b. Rows that had an explicit height before but must go back to the
#     sheet default (no ht/customHeight attribute at all)
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(23).AutoFit()

# ---- 5. Column widths ----
# Column A & B currently share one <col min="1" max="2" .../> entry;
# touching column B's width forces the sheet to record column A on
# its own (min="1" max="1"), matching the target layout, while column
# B keeps its original 60.7109375 width.
$ws.Columns.Item(2).ColumnWidth = 60.7109375

